$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.406.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.42%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.106.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  -0.93%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "343.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("E6").Value = "  -0.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5323"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4443"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.90"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09413"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.173"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.574"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.927"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.076.63"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "101.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.34%  "
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06689"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.334"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.420.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.316"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.821"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.31%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.528"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.148"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.675"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.03%  "
$ws.Range("E33").Value = "  +0.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.270"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.848"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02657"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.36%  "
$ws.Range("E38").Value = "  +1.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.74"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.7030"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.349"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2227"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6877"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.350"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.68%  "
$ws.Range("E46").Value = "  -0.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.382"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +18.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.640"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000350"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.224"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.58%  "
$ws.Range("E51").Value = "  +0.23%  "
